# ---------------------------------------------------------------------------
# Adds a "2022-Q3" quarterly sheet to the workbook:
#   1. A brand-new worksheet "2022-Q3" is inserted right after "总计" and
#      before "2022-Q2" (all later sheets shift right by one tab, unchanged
#      in content).
#   2. The "总计" (summary) sheet gets a new row for 2022-Q3 at the top of
#      the data block, every other quarter's row shifts down by one, and a
#      duplicate trailing row for "2020-Q4" is appended so the last original
#      row is preserved once the data has shifted.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q3" worksheet before the existing "2022-Q2" sheet.
# ---------------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($refSheet)
$q3.Name = "2022-Q3"

# Header row.
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
$q3.Range("B1:H1").NumberFormat = "@"
for ($c = 0; $c -lt $headers.Count; $c++) {
    $q3.Cells.Item(1, 2 + $c).Value = $headers[$c]
}
$q3.Range("B1:H1").Style = "Normal"

# Fund holdings data (columns: A idx, B code, C name, D size, E stock-pos,
# F pos-share, G market-value, H pos-rank).
$rows = @(
    @("0","161611","融通内需驱动混合A/B","8.72","90.68","4.87","0.4247","3"),
    @("1","014772","中泰红利价值一年持有混合","5.05","93.55","4.99","0.2520","8"),
    @("2","014771","中泰红利优选一年持有混合","5.02","93.69","4.92","0.2470","9"),
    @("3","014109","融通内需驱动混合C","4.06","90.68","4.87","0.1977","3"),
    @("4","002252","融通成长30灵活配置混合A/B","1.56","93.13","5.09","0.0794","2"),
    @("5","014106","融通成长30灵活配置混合C","1.34","93.13","5.09","0.0682","2"),
    @("6","008115","天弘中证红利低波动100指数C","2.44","94.56","1.66","0.0405","10"),
    @("7","008114","天弘中证红利低波动100指数A","1.89","94.56","1.66","0.0314","10"),
    @("8","515100","景顺长城中证红利低波动100ETF","1.62","98.63","1.74","0.0282","10"),
    @("9","013611","工银民瑞一年持有混合A","3.05","21.97","0.82","0.0250","10"),
    @("10","512590","浦银安盛中证高股息精选ETF","0.45","90.87","1.95","0.0088","7"),
    @("11","164811","工银瑞信中证京津冀协同发展主题指数（LOF）A","0.12","93.09","3.11","0.0037","5"),
    @("12","164825","工银瑞信中证京津冀协同发展主题指数（LOF）C","0.03","93.09","3.11","0.0009","5"),
    @("13","013612","工银民瑞一年持有混合C","0.10","21.97","0.82","0.0008","10")
)

$lastRow = 1 + $rows.Count
$q3.Range("B2:G$lastRow").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $q3.Cells.Item($r, 1).Value = [int]$row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = [int]$row[7]
}

$q3.Range("B2:G$lastRow").Style = "Normal"

# ---------------------------------------------------------------------------
# 2) Rewrite the "总计" summary sheet with the new quarter prepended.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$summary = @(
    @("2022-Q3","14","1.41"),
    @("2022-Q2","6","0.86"),
    @("2022-Q1","19","2.35"),
    @("2021-Q4","3","0.16"),
    @("2021-Q3","5","0.16"),
    @("2021-Q2","3","0.01"),
    @("2021-Q1","7","0.09"),
    @("2020-Q4","7","0.09")
)

# Copy the formatting of the existing last row down to the brand-new row 9
# so the new "A" index cell picks up the same bold/border style as A2:A8.
$total.Range("A8").Copy($total.Range("A9"))
$total.Range("B8:D8").Copy($total.Range("B9:D9"))

for ($i = 0; $i -lt $summary.Count; $i++) {
    $r = $i + 2
    $item = $summary[$i]
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $item[0]
    $total.Cells.Item($r, 3).Value = [int]$item[1]
    $total.Cells.Item($r, 4).Value = [double]$item[2]
}

# ---------------------------------------------------------------------------
# Keep "总计" as the active sheet, matching the original selection.
# ---------------------------------------------------------------------------
$total.Activate()
